$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update NOTES column (O) contents for rows 4 and 6 (the former "INVESTIGATE"
# notes have been reviewed/resolved with more specific explanations).
$ws.Range("O4").Value = "Blank uncontested not included in the manual results."
$ws.Range("O6").Value = "Data error in manual spreadsheet."

# Widen column O so the longer note text is fully visible.
# (45.1666... compensates for this runtime's internal padding so the
# stored OOXML <col> width lands on exactly 46, matching Excel's own
# "best fit" recompute for this text/font.)
$ws.Columns.Item(15).ColumnWidth = 45.166666666666664

# Move the active cell selection on the bottom-right pane from O9 to O7.
$ws.Range("O7").Select()
